$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.90"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.90"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.380"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05894"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.384"
$ws.Range("E6").Value = "5GateTokenGT"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.503"
$ws.Range("E7").Value = "6KuCoinTokenKCS"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8105"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9285"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1417"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07415"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03048"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03054"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.890"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001562"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04694"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005976"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005907"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001242"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004740"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00008796"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.555"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3229"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1329"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002651"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03891"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006317"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1070"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002799"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008557"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005208"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7484"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001943"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
